$wb = $excel.ActiveWorkbook

# --- Sheet "produtos": clear C92 (now empty) and append new rows 93-97 ---
$ws1 = $wb.Worksheets.Item("produtos")

# C92 becomes a truly empty cell (no inlineStr placeholder) instead of an empty inline string
$ws1.Range("C92").Value = ""

$produtos = @(
    @(92, "rr", "kg"),
    @(93, "alcool isoproponiu", "5 L"),
    @(94, "alcool iso  do bom", "PCT C/10"),
    @(95, "arroz soltinho", "kg"),
    @(96, "feijão pretola", "kg")
)

$r = 93
foreach ($item in $produtos) {
    $ws1.Cells.Item($r, 1).Value = $item[0]
    $ws1.Cells.Item($r, 2).Value = $item[1]
    $ws1.Cells.Item($r, 4).Value = $item[2]
    $ws1.Cells.Item($r, 5).Value = 0
    $r++
}

# --- Sheet "movimentos": append new rows 15-16 ---
$ws2 = $wb.Worksheets.Item("movimentos")

$movimentos = @(
    @(14, 2, "ENTRADA", 11, "2025-12-16 16:39:59"),
    @(15, 4, "ENTRADA", 11, "2025-12-16 16:39:59")
)

$r = 15
foreach ($item in $movimentos) {
    $ws2.Cells.Item($r, 1).Value = $item[0]
    $ws2.Cells.Item($r, 2).Value = $item[1]
    $ws2.Cells.Item($r, 3).Value = $item[2]
    $ws2.Cells.Item($r, 4).Value = $item[3]
    $ws2.Cells.Item($r, 5).Value = $item[4]
    $r++
}
